$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price/Volume columns to text so numeric-looking strings
# (e.g. "1.00", "67.158.15") are preserved exactly as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '67.074.08'
$ws.Range("E2").Value = '  +8.38%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.517.49'
$ws.Range("E3").Value = '  +12.08%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '190.21'
$ws.Range("E5").Value = '  +13.16%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '550.57'
$ws.Range("E6").Value = '  +6.93%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.509.02'
$ws.Range("E7").Value = '  +11.70%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.606'
$ws.Range("E8").Value = '  +3.92%  '

$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.632'
$ws.Range("E10").Value = '  +7.21%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.151'
$ws.Range("E11").Value = '  +19.12%  '

$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '54.89'
$ws.Range("E12").Value = '  +6.45%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000269'
$ws.Range("E13").Value = '  +9.72%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '9.37'
$ws.Range("E14").Value = '  +5.98%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.071.26'
$ws.Range("E15").Value = '  +12.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.519.30'
$ws.Range("E16").Value = '  +12.64%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.121'
$ws.Range("E17").Value = '  +6.14%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '67.023.93'
$ws.Range("E18").Value = '  +8.93%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '18.16'
$ws.Range("E19").Value = '  +8.05%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '11.87'
$ws.Range("E20").Value = '  +9.94%  '

$ws.Range("B21").Value = 'Polygon'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D21").Value = '0.992'
$ws.Range("E21").Value = '  +4.48%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '425.43'
$ws.Range("E22").Value = '  +18.79%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '84.99'
$ws.Range("E23").Value = '  +7.35%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '3.91'
$ws.Range("E24").Value = '  +6.99%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '4.19'
$ws.Range("E25").Value = '  +9.03%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '11.12'
$ws.Range("E26").Value = '  +0.25%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '2.91'
$ws.Range("E27").Value = '  +13.34%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '11.94'
$ws.Range("E28").Value = '  +8.79%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '8.87'
$ws.Range("E29").Value = '  +11.49%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '30.18'
$ws.Range("E30").Value = '  +8.86%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '653.93'
$ws.Range("E31").Value = '  +2.49%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '6.65'
$ws.Range("E32").Value = '  +5.83%  '

$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '11.72'
$ws.Range("E33").Value = '  +5.55%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  +8.10%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '59.46'
$ws.Range("E35").Value = '  +5.27%  '

$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '38.53'
$ws.Range("E36").Value = '  +8.16%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0817'
$ws.Range("E37").Value = '  +20.55%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = '0.389'
$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.139'
$ws.Range("E40").Value = '  +14.96%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '3.28'
$ws.Range("E41").Value = '  +15.56%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.017.83'
$ws.Range("E43").Value = '  +6.37%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.64'
$ws.Range("E44").Value = '  +4.62%  '

$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").Value = '  +15.87%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.27'
$ws.Range("E46").Value = '  +11.37%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0416'
$ws.Range("E47").Value = '  +9.39%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.71'
$ws.Range("E48").Value = '  +4.05%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.130'
$ws.Range("E49").Value = '  +7.75%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '8.74'
$ws.Range("E50").Value = '  +18.15%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '140.68'
$ws.Range("E51").Value = '  +8.02%  '

# Restore the original (default) style on the Price/Volume columns so the
# temporary text-number-format does not leave a stray style behind.
$ws.Range("D2:E51").Style = "Normal"
